# edit.ps1 -- apply the "Add files via upload" revision to the ODIM
# architecture slide:
#   * rebrand four small component labels (drop the "HPE"/"Aruba" vendor
#     prefixes) and resize/reposition their textboxes to match the
#     shorter text,
#   * drop the "HPE" prefix from the big product-name textbox (leaving
#     the trailing "(TM)" run, its formatting, and the rest of the
#     paragraph untouched).
#
# Shape lookup is done by the PowerPoint shape Id (stable identifier;
# survives group re-indexing) via a small recursive helper, since several
# of the target shapes live inside a nested group ("Group 4").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeById {
    param($shapes, [int]$id)
    $count = $shapes.Count
    foreach ($candidate in $shapes) {
        if ($candidate.Id -eq $id) {
            return $candidate
        }
        if ($candidate.Type -eq 6) {
            # msoGroup -- recurse into its GroupItems
            $found = Find-ShapeById $candidate.GroupItems $id
            if ($found) {
                return $found
            }
        }
    }
    return $null
}

# --- "HPE Servers" -> "Servers" (TextBox 20, id 21) ------------------------
$sh = Find-ShapeById $s.Shapes 21
$sh.TextFrame.TextRange.Text = "Servers"
$sh.Left = 137.17536173070823
$sh.Top = 404.77259842519686
$sh.Width = 50.553858267716535
$sh.Height = 21.810944881889764

# --- "Aruba Fabric" -> "Fabric" (TextBox 21, id 22) ------------------------
$sh = Find-ShapeById $s.Shapes 22
$sh.TextFrame.TextRange.Text = "Fabric"
$sh.Left = 234.79102362204725
$sh.Top = 404.77259842519686
$sh.Width = 44.0055905511811
$sh.Height = 21.810944881889764

# --- "HPE iLO Plugin" -> "BMC Plugin" (TextBox 39, id 40) ------------------
$sh = Find-ShapeById $s.Shapes 40
$sh.TextFrame.TextRange.Text = "BMC Plugin"
$sh.Left = 134.71606299212598
$sh.Top = 300.58314960629923
$sh.Width = 61.4948031496063
$sh.Height = 25.446063092125986

# --- "Aruba Plugin" -> "Fabric Plugin" (TextBox 40, id 41) -----------------
$sh = Find-ShapeById $s.Shapes 41
$sh.TextFrame.TextRange.Text = "Fabric Plugin"
$sh.Left = 219.72622047244096
$sh.Top = 300.58314960629923
$sh.Width = 66.79606299212598
$sh.Height = 25.446063092125986

# --- "HPE Resource Aggregator for ODIM(TM)" -> "Resource Aggregator for ODIM(TM)"
# (TextBox 76, id 77) -- only the first run's text changes; the trailing
# trademark-symbol run (superscript, separately colored) is left alone, so
# address just the first 32 characters of the range.
$sh = Find-ShapeById $s.Shapes 77
$tr = $sh.TextFrame.TextRange
$firstRun = $tr.Characters(1, 32)
$firstRun.Text = "Resource Aggregator for ODIM"
